$d = $word.ActiveDocument

# Locate the end of the "LOQ4057: ..." requirement line (the paragraph mark
# that follows it is where the text to remove begins).
$rStart = $d.Content
$rStart.Find.Execute(
    "LOQ4057: Operações Unitárias III (Requisito fraco)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteStart = $rStart.End + 1

# Locate the end of the trailing copyright/footer line; its paragraph mark is
# where the text to remove ends (everything up through & including that mark
# goes away, along with the "Ver no Jupiter..." paragraph in between).
$rEnd = $d.Content
$rEnd.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$deleteEnd = $rEnd.End + 1

$d.Range($deleteStart, $deleteEnd).Delete()
